# Tracking.xlsx update: "Update task tracking va task trong tet"
#  - Attendent check: add a new attendance row (23/01/04) below 21/01/14,
#    and fill in the status cells for 21/01/14 that were left blank.
#  - Deadline: add a new deadline row (23/01/14), copied from the existing
#    19/01/14 deadline row.
#  - Leave the workbook with the "Deadline" sheet active/selected.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Attendent check")
$ws2 = $wb.Worksheets.Item("Deadline")

# PasteSpecial modes used below (standard Excel xlPasteType values):
#   -4122 = xlPasteFormats   (formats only, keep existing/blank value)
#   -4163 = xlPasteValues    (values only)

# ---- Sheet "Attendent check": finish formatting row 10 (21/01/14) ----
$ws1.Range("B9").Copy() | Out-Null
$ws1.Range("B10").PasteSpecial(-4122) | Out-Null
$ws1.Range("D6").Copy() | Out-Null
$ws1.Range("C10").PasteSpecial(-4122) | Out-Null
$ws1.Range("B9").Copy() | Out-Null
$ws1.Range("D10").PasteSpecial(-4122) | Out-Null
$ws1.Range("B9").Copy() | Out-Null
$ws1.Range("E10").PasteSpecial(-4122) | Out-Null

# ---- Sheet "Attendent check": new row 11 for 23/01/04 ----
$ws1.Range("A11").Value = "23/01/04"
$ws1.Range("B9").Copy() | Out-Null
$ws1.Range("B11").PasteSpecial(-4122) | Out-Null
$ws1.Range("B9").Copy() | Out-Null
$ws1.Range("C11").PasteSpecial(-4122) | Out-Null
$ws1.Range("D9").Copy() | Out-Null
$ws1.Range("D11").PasteSpecial(-4122) | Out-Null
$ws1.Range("D6").Copy() | Out-Null
$ws1.Range("E11").PasteSpecial(-4122) | Out-Null

# ---- Sheet "Deadline": new row 10 for 23/01/14, mirroring row 9 ----
$ws2.Range("A9:E9").Copy() | Out-Null
$ws2.Range("A10:E10").PasteSpecial(-4163) | Out-Null
$ws2.Range("A9:E9").Copy() | Out-Null
$ws2.Range("A10:E10").PasteSpecial(-4122) | Out-Null
$ws2.Range("A10").Value = "23/01/14"

# ---- Selection / active sheet bookkeeping ----
$ws1.Range("E11").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("E10").Select() | Out-Null
